# RSTK-9618 - Derived Component Receipt template changes
# Renames "Pro-Disassembley (Lot and serial track)" -> "Pro-Disassembley1 (Lot and serial track)"
# and removes the now-obsolete "Pro-Disassembley Serial (Serial track)" rows/variant across all
# three worksheets, plus refreshes the active-sheet/selection bookmarks left behind by Excel.

$wb = $excel.ActiveWorkbook

$wsCreate    = $wb.Worksheets.Item("Create Disassembly WO")
$wsDerived   = $wb.Worksheets.Item("Derived Components")
$wsConsum    = $wb.Worksheets.Item("Consumable Components")

# ---------------------------------------------------------------------------
# Sheet 1 - "Create Disassembly WO"
# ---------------------------------------------------------------------------
$wsCreate.Range("B2").Value = "Pro-Disassembley1 (Lot and serial track)"
$wsCreate.Columns.Item(2).ColumnWidth = 33.14

# ---------------------------------------------------------------------------
# Sheet 2 - "Derived Components"
# ---------------------------------------------------------------------------
# Update the remaining "lot and serial" rows (2 & 3) to point at the renamed variant
# and bump their quantities from 2 -> 3.
$wsDerived.Range("A2").Value = "Pro-Disassembley1 (Lot and serial track)"
$wsDerived.Range("C2").Value = 3
$wsDerived.Range("D2").Value = 3

$wsDerived.Range("A3").Value = "Pro-Disassembley1 (Lot and serial track)"
$wsDerived.Range("C3").Value = 3
$wsDerived.Range("D3").Value = 3

# Remove the obsolete "Pro-Disassembley Serial (Serial track)" rows (4 & 5), keeping only
# the style-carrying, now-empty F/J cells behind.
$wsDerived.Range("A4").ClearContents()
$wsDerived.Range("B4").ClearContents()
$wsDerived.Range("C4").ClearContents()
$wsDerived.Range("D4").ClearContents()
$wsDerived.Range("F4").ClearContents()

$wsDerived.Range("A5").ClearContents()
$wsDerived.Range("B5").ClearContents()
$wsDerived.Range("C5").ClearContents()
$wsDerived.Range("D5").ClearContents()
$wsDerived.Range("F5").ClearContents()
$wsDerived.Range("I5").ClearContents()
$wsDerived.Range("J5").ClearContents()

$wsDerived.Columns.Item(1).ColumnWidth = 33.14

# ---------------------------------------------------------------------------
# Sheet 3 - "Consumable Components"
# ---------------------------------------------------------------------------
# Rows 2-4: point at the renamed variant and fill in the explicit boolean/zero values that
# Excel now writes for the previously-blank cells.
$wsConsum.Range("A2").Value = "Pro-Disassembley1 (Lot and serial track)"
$wsConsum.Range("D2").Value = $false
$wsConsum.Range("E2").Value = $false
$wsConsum.Range("G2").Value = 0
$wsConsum.Range("H2").Value = 0

$wsConsum.Range("A3").Value = "Pro-Disassembley1 (Lot and serial track)"
$wsConsum.Range("C3").Value = $false
$wsConsum.Range("E3").Value = $false

$wsConsum.Range("A4").Value = "Pro-Disassembley1 (Lot and serial track)"
$wsConsum.Range("C4").Value = $false
$wsConsum.Range("D4").Value = $false

# Remove the obsolete "Pro-Disassembley Serial (Serial track)" rows (5, 6 & 7) entirely.
$wsConsum.Rows.Item(5).Delete()
$wsConsum.Rows.Item(5).Delete()
$wsConsum.Rows.Item(5).Delete()

$wsConsum.Columns.Item(1).ColumnWidth = 33.14

# ---------------------------------------------------------------------------
# Selections / active sheet bookkeeping (mirrors what Excel persists on save)
# ---------------------------------------------------------------------------
$wsDerived.Range("A3").Select()
$wsCreate.Range("B2").Select()
$wsConsum.Range("D16").Select()
$wsConsum.Select()
